$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("quality_comparison")
$ws2 = $wb.Worksheets.Item("computational_comparison")

# ---------------------------------------------------------------------
# Build the two new header-border looks once (on quality_comparison!C1/D1)
# and stamp them onto every other cell that needs them by copying the
# already-built format. C1 gets a top+bottom rule, D1 (right edge of the
# merged B:D header band) gets top+bottom+right.
# ---------------------------------------------------------------------
$ws1.Range("C1").Style = "Normal"
$ws1.Range("C1").Borders(8).LineStyle = 1
$ws1.Range("C1").Borders(9).LineStyle = 1

$ws1.Range("D1").Style = "Normal"
$ws1.Range("D1").Borders(10).LineStyle = 1
$ws1.Range("D1").Borders(9).LineStyle = 1
$ws1.Range("D1").Borders(8).LineStyle = 1

# computational_comparison has two header bands (B:D and E:G) that both
# need the same treatment.
$ws1.Range("C1").Copy()
$ws2.Range("C1").PasteSpecial(-4122)

$ws1.Range("D1").Copy()
$ws2.Range("D1").PasteSpecial(-4122)

$ws1.Range("C1").Copy()
$ws2.Range("F1").PasteSpecial(-4122)

$ws1.Range("D1").Copy()
$ws2.Range("G1").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# Anonymize the "fedcore" column headers on both sheets.
# ---------------------------------------------------------------------
$ws1.Range("C2").Value = "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# ---------------------------------------------------------------------
# G5 was an empty placeholder cell (no value, no style) — drop it.
# ---------------------------------------------------------------------
$ws2.Range("G5").ClearContents()
